# Regenerate save_data to use K instead of Strike#, recalculated values
# written into column G (header "K") for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 2
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 4
$ws.Range("G9").Value = 3
